# Hello Patient Tutorial - Navigation added
# Adds three "section label" textboxes to the right of the XML/JSON example
# boxes on slide 2: "Resource identity and metadata", "Human readable
# summary" and "Standard data items".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# EMU -> point conversion (PowerPoint COM Shape geometry is expressed in
# points; 1 pt == 12700 EMU).
$EMU = 12700.0

# The slide already has shape ids 1,4,5,6,8,9,10,11,15,16,17,18 in use.
# PowerPoint hands out the next *unused* id (starting at 2) to every shape
# Shapes.Add* call creates, regardless of deletions, so to land on ids
# 19/20/21 (matching the authored deck) burn through the lower unused ids
# (2,3,7,12,13,14) first with throw-away textboxes, then delete them.
$scratch = @()
for ($i = 0; $i -lt 6; $i++) {
    $scratch += $s.Shapes.AddTextbox(1, 0, 0, 1, 1)
}

# NOTE: this interpreter only binds function parameters positionally, so
# avoid `-Name value` style calls here.
function Add-SectionLabel {
    param(
        [string]$Name,
        [double]$Left,
        [double]$Top,
        [double]$Width,
        [double]$Height,
        [string]$Text
    )

    $shp = $s.Shapes.AddTextbox(1, $Left / $EMU, $Top / $EMU, $Width / $EMU, $Height / $EMU)
    $shp.Name = $Name
    $shp.TextFrame.WordWrap = $false
    $shp.TextFrame.AutoSize = 1
    $shp.TextFrame.TextRange.Text = $Text
    $shp.TextFrame.TextRange.Font.Size = 14
    $shp.Fill.Visible = $false
    return $shp
}

$t1 = Add-SectionLabel "TextBox 18" 6168305 2106209 2520917 307777 "Resource identity and metadata"
$t2 = Add-SectionLabel "TextBox 19" 6168305 2913349 2121093 307777 "Human readable summary"
$t3 = Add-SectionLabel "TextBox 20" 6168305 4191273 1656811 307777 "Standard data items"

foreach ($scr in $scratch) {
    $scr.Delete()
}

Write-Output "Added $($t1.Name) (id=$($t1.Id)), $($t2.Name) (id=$($t2.Id)), $($t3.Name) (id=$($t3.Id)); slide now has $($s.Shapes.Count) shapes"
